$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for the new data block so numeric-looking values
# (including SKUs with leading zeros) are stored as text, matching the source rows.
$ws.Range("A3:E15").NumberFormat = "@"

$ws.Cells.Item(3, 1).Value = "462130"
$ws.Cells.Item(3, 2).Value = "Joe Tea - Peach"
$ws.Cells.Item(3, 3).Value = "4"
$ws.Cells.Item(3, 4).Value = "22.49"
$ws.Cells.Item(3, 5).Value = "89.96"

$ws.Cells.Item(4, 1).Value = "462225"
$ws.Cells.Item(4, 2).Value = "Joe Tea - Mango Lemonade"
$ws.Cells.Item(4, 3).Value = "2"
$ws.Cells.Item(4, 4).Value = "22.49"
$ws.Cells.Item(4, 5).Value = "44.98"

$ws.Cells.Item(5, 1).Value = "462120"
$ws.Cells.Item(5, 2).Value = "Joe Tea - Lemon"
$ws.Cells.Item(5, 3).Value = "3"
$ws.Cells.Item(5, 4).Value = "22.49"
$ws.Cells.Item(5, 5).Value = "67.47"

$ws.Cells.Item(6, 1).Value = "462180"
$ws.Cells.Item(6, 2).Value = "Joe Tea - Kiwi Strawberry"
$ws.Cells.Item(6, 3).Value = "4"
$ws.Cells.Item(6, 4).Value = "22.49"
$ws.Cells.Item(6, 5).Value = "89.96"

$ws.Cells.Item(7, 1).Value = "462115"
$ws.Cells.Item(7, 2).Value = "Joe Tea - Half & Half"
$ws.Cells.Item(7, 3).Value = "5"
$ws.Cells.Item(7, 4).Value = "22.49"
$ws.Cells.Item(7, 5).Value = "112.45"

$ws.Cells.Item(8, 1).Value = "462110"
$ws.Cells.Item(8, 2).Value = "Joe Tea - Ginseng Green"
$ws.Cells.Item(8, 3).Value = "1"
$ws.Cells.Item(8, 4).Value = "22.49"
$ws.Cells.Item(8, 5).Value = "22.49"

$ws.Cells.Item(9, 1).Value = "462105"
$ws.Cells.Item(9, 2).Value = "Joe Tea - Classic Lemonade"
$ws.Cells.Item(9, 3).Value = "5"
$ws.Cells.Item(9, 4).Value = "22.49"
$ws.Cells.Item(9, 5).Value = "112.45"

$ws.Cells.Item(10, 1).Value = "462100"
$ws.Cells.Item(10, 2).Value = "Joe Tea - Black Unsweetened"
$ws.Cells.Item(10, 3).Value = "4"
$ws.Cells.Item(10, 4).Value = "22.49"
$ws.Cells.Item(10, 5).Value = "89.96"

$ws.Cells.Item(11, 1).Value = "462175"
$ws.Cells.Item(11, 2).Value = "Joe Tea - Black Cherry"
$ws.Cells.Item(11, 3).Value = "12"
$ws.Cells.Item(11, 4).Value = "22.49"
$ws.Cells.Item(11, 5).Value = "269.88"

$ws.Cells.Item(12, 1).Value = "456501"
$ws.Cells.Item(12, 2).Value = "Bragg - Honey & Green Tea"
$ws.Cells.Item(12, 3).Value = "1"
$ws.Cells.Item(12, 4).Value = "27.74"
$ws.Cells.Item(12, 5).Value = "27.74"

$ws.Cells.Item(13, 1).Value = "456505"
$ws.Cells.Item(13, 2).Value = "Bragg - Ginger Lemon Honey"
$ws.Cells.Item(13, 3).Value = "1"
$ws.Cells.Item(13, 4).Value = "30.99"
$ws.Cells.Item(13, 5).Value = "30.99"

$ws.Cells.Item(14, 1).Value = "053365"
$ws.Cells.Item(14, 2).Value = "DV - Yogurt Mini Pretzel"
$ws.Cells.Item(14, 3).Value = "1"
$ws.Cells.Item(14, 4).Value = "29.63"
$ws.Cells.Item(14, 5).Value = "29.63"

$ws.Cells.Item(15, 1).Value = "456090"
$ws.Cells.Item(15, 2).Value = "Employee Water"
$ws.Cells.Item(15, 3).Value = "15"
$ws.Cells.Item(15, 4).Value = "2.93"
$ws.Cells.Item(15, 5).Value = "43.95"
